# atualizacao 16 nov 2020
# Appends the newest daily/monthly readings to both sheets:
#  - "Mensal": one new row (A14:D14) for the 2020-11 monthly summary
#  - "Diario": fifteen new rows (A368:D382) for 2020-11-01 .. 2020-11-15

$wb = $excel.ActiveWorkbook

# ---- Mensal sheet: add row 14 ----
$wsMensal = $wb.Worksheets.Item("Mensal")
$lastMensalRow = 13
$newMensalRow = 14

# Copy the formatting (incl. the date number format style) from the last
# existing data row down onto the new row before writing values into it.
$wsMensal.Range("A" + $lastMensalRow).Copy()
$wsMensal.Range("A" + $newMensalRow).PasteSpecial(-4122)

$wsMensal.Cells.Item($newMensalRow, 1).Value = 44150
$wsMensal.Cells.Item($newMensalRow, 2).Value = 39.65
$wsMensal.Cells.Item($newMensalRow, 3).Value = 133.23
$wsMensal.Cells.Item($newMensalRow, 4).Value = -70.23999999999999

# ---- Diario sheet: add rows 368..382 ----
$wsDiario = $wb.Worksheets.Item("Diario")
$lastDiarioRow = 367

$newDiarioRows = @(
    @(44136, 36.88, 133.23, -72.31999999999999),
    @(44137, 41.61, 133.23, -68.77),
    @(44138, 41.45, 133.23, -68.89),
    @(44139, 43.82, 133.23, -67.11),
    @(44140, 41.96, 133.23, -68.51000000000001),
    @(44141, 34.26, 133.23, -74.29000000000001),
    @(44142, 29.21, 133.23, -78.08),
    @(44143, 25.39, 133.23, -80.95),
    @(44144, 22.3, 133.23, -83.26000000000001),
    @(44145, 22.71, 133.23, -82.95999999999999),
    @(44146, 25.55, 133.23, -80.81999999999999),
    @(44147, 51.4, 133.23, -61.42),
    @(44148, 68.38, 133.23, -48.68),
    @(44149, 59.78, 133.23, -55.13),
    @(44150, 50.12, 133.23, -62.38)
)

$row = $lastDiarioRow + 1
foreach ($values in $newDiarioRows) {
    $wsDiario.Range("A" + $lastDiarioRow).Copy()
    $wsDiario.Range("A" + $row).PasteSpecial(-4122)

    $wsDiario.Cells.Item($row, 1).Value = $values[0]
    $wsDiario.Cells.Item($row, 2).Value = $values[1]
    $wsDiario.Cells.Item($row, 3).Value = $values[2]
    $wsDiario.Cells.Item($row, 4).Value = $values[3]

    $row = $row + 1
}
